# Adding rules to engagement creation form
# Inserts two new rule rows ("Création imputation" / "Création apurement")
# into the "TODO" sheet, right after the existing "Création engagement" rules
# (rows 5-6), and re-colors those existing rules to match the sheet's
# highlight (green) style while the new rows take over the previous
# "no fill" style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")

# Insert two blank rows above the current row 7 ("Apurements" section),
# shifting everything below down by two rows (old row 7 -> new row 9, etc.)
$ws.Rows.Item(7).Insert()
$ws.Rows.Item(7).Insert()

# Fill in the two new rows. B8 is set first so that "Création apurement"
# is registered as a shared string before "Création imputation", matching
# the order the two new entries should appear in the workbook.
$ws.Range("B8").Value = "Création apurement"
$ws.Range("B7").Value = "Création imputation"

# The newly inserted rows inherit the formatting from the row above them
# (the "no fill" look previously used by rows 5-6), which is what we want
# for these new rows. Now re-color the original "Création engagement"
# rows (B5:C6) with the sheet's green highlight fill.
$ws.Range("B5:C6").Interior.Color = 5296274

# Update the sheet's active selection to B14 (the new location of the
# "Gestion des autorisation au niveau du backend" row).
$ws.Range("B14").Select()
